# issue #5: stock data from json to db
# Add "category", "source_file", and "index" columns to the 股票 (stock) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# Insert a new column I for "category" (between property_category @H and date @I).
$ws.Columns.Item(9).Insert()

# After this insert, legislator_id moved from K to L.
# Insert two new columns after it (M = source_file, N = index).
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(13).Insert()

# --- Header row ---
$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data rows ---
# Row 2 (index 63)
$ws.Range("I2").Value = "normal"
$ws.Range("M2").Value = "tmp59331"
$ws.Range("N2").Value = 63

# Row 3 (index 64)
$ws.Range("I3").Value = "normal"
$ws.Range("M3").Value = "tmp59331"
$ws.Range("N3").Value = 64

# Row 4 (index 65)
$ws.Range("I4").Value = "normal"
$ws.Range("M4").Value = "tmp59331"
$ws.Range("N4").Value = 65

# Row 5 (index 66)
$ws.Range("I5").Value = "normal"
$ws.Range("M5").Value = "tmp59331"
$ws.Range("N5").Value = 66

# Row 6 (index 67)
$ws.Range("I6").Value = "normal"
$ws.Range("M6").Value = "tmp59331"
$ws.Range("N6").Value = 67

# Row 7 (index 68)
$ws.Range("I7").Value = "normal"
$ws.Range("M7").Value = "tmp59331"
$ws.Range("N7").Value = 68

Write-Output "done"
